$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cell values (sending/target cluster reassignments and recomputed TPM-based statistics)
$ws.Cells.Item(2, 7).Value2 = 2.060111333333333
$ws.Cells.Item(2, 8).Value2 = 6.180334
$ws.Cells.Item(2, 9).Value2 = 0.2095457297481522
$ws.Cells.Item(2, 10).Value2 = 0.2095457297481522
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 12).Value2 = 1
$ws.Cells.Item(2, 13).Value2 = 34.97976933333334
$ws.Cells.Item(2, 14).Value2 = 104.939308
$ws.Cells.Item(2, 15).Value2 = 0.4352965780925344
$ws.Cells.Item(2, 16).Value2 = 0.4352965780925344
$ws.Cells.Item(2, 17).Value2 = 72.06221924098578
$ws.Cells.Item(2, 18).Value2 = 648.5599731688721
$ws.Cells.Item(2, 19).Value2 = 0.09121453911327362
$ws.Cells.Item(2, 20).Value2 = 0.09121453911327364
$ws.Cells.Item(3, 7).Value2 = 2.060111333333333
$ws.Cells.Item(3, 8).Value2 = 6.180334
$ws.Cells.Item(3, 9).Value2 = 0.2095457297481522
$ws.Cells.Item(3, 10).Value2 = 0.2095457297481522
$ws.Cells.Item(3, 14).Value2 = 61.03014900000001
$ws.Cells.Item(3, 15).Value2 = 0.2531579017099818
$ws.Cells.Item(3, 16).Value2 = 0.2531579017099818
$ws.Cells.Item(3, 17).Value2 = 41.90963387664068
$ws.Cells.Item(3, 18).Value2 = 377.1867048897661
$ws.Cells.Item(3, 19).Value2 = 0.05304815725532912
$ws.Cells.Item(3, 20).Value2 = 0.05304815725532912
$ws.Cells.Item(4, 7).Value2 = 2.060111333333333
$ws.Cells.Item(4, 8).Value2 = 6.180334
$ws.Cells.Item(4, 9).Value2 = 0.2095457297481522
$ws.Cells.Item(4, 10).Value2 = 0.2095457297481522
$ws.Cells.Item(4, 13).Value2 = 25.035323
$ws.Cells.Item(4, 14).Value2 = 75.105969
$ws.Cells.Item(4, 15).Value2 = 0.3115455201974837
$ws.Cells.Item(4, 16).Value2 = 0.3115455201974837
$ws.Cells.Item(4, 17).Value2 = 51.57555264596067
$ws.Cells.Item(4, 18).Value2 = 464.179973813646
$ws.Cells.Item(4, 19).Value2 = 0.06528303337954941
$ws.Cells.Item(4, 20).Value2 = 0.06528303337954941
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 7).Value2 = 4.24427
$ws.Cells.Item(5, 8).Value2 = 12.73281
$ws.Cells.Item(5, 9).Value2 = 0.431709024657012
$ws.Cells.Item(5, 10).Value2 = 0.431709024657012
$ws.Cells.Item(5, 11).Value2 = 3
$ws.Cells.Item(5, 12).Value2 = 1
$ws.Cells.Item(5, 13).Value2 = 34.97976933333334
$ws.Cells.Item(5, 14).Value2 = 104.939308
$ws.Cells.Item(5, 15).Value2 = 0.4352965780925344
$ws.Cells.Item(5, 16).Value2 = 0.4352965780925344
$ws.Cells.Item(5, 17).Value2 = 148.4635855883867
$ws.Cells.Item(5, 18).Value2 = 1336.17227029548
$ws.Cells.Item(5, 19).Value2 = 0.1879214611648629
$ws.Cells.Item(5, 20).Value2 = 0.1879214611648629
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 9).Value2 = 0.431709024657012
$ws.Cells.Item(6, 10).Value2 = 0.431709024657012
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 12).Value2 = 1
$ws.Cells.Item(6, 13).Value2 = 20.343383
$ws.Cells.Item(6, 14).Value2 = 61.03014900000001
$ws.Cells.Item(6, 15).Value2 = 0.2531579017099818
$ws.Cells.Item(6, 16).Value2 = 0.2531579017099818
$ws.Cells.Item(6, 17).Value2 = 86.34281016541001
$ws.Cells.Item(6, 18).Value2 = 777.0852914886901
$ws.Cells.Item(6, 19).Value2 = 0.109290550831432
$ws.Cells.Item(6, 20).Value2 = 0.109290550831432
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 9).Value2 = 0.431709024657012
$ws.Cells.Item(7, 10).Value2 = 0.431709024657012
$ws.Cells.Item(7, 13).Value2 = 25.035323
$ws.Cells.Item(7, 14).Value2 = 75.105969
$ws.Cells.Item(7, 15).Value2 = 0.3115455201974837
$ws.Cells.Item(7, 16).Value2 = 0.3115455201974837
$ws.Cells.Item(7, 17).Value2 = 106.25667034921
$ws.Cells.Item(7, 18).Value2 = 956.3100331428901
$ws.Cells.Item(7, 19).Value2 = 0.1344970126607171
$ws.Cells.Item(7, 20).Value2 = 0.1344970126607171
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 7).Value2 = 2.602283
$ws.Cells.Item(8, 8).Value2 = 7.806849
$ws.Cells.Item(8, 9).Value2 = 0.2646931170287289
$ws.Cells.Item(8, 10).Value2 = 0.2646931170287289
$ws.Cells.Item(8, 13).Value2 = 34.97976933333334
$ws.Cells.Item(8, 14).Value2 = 104.939308
$ws.Cells.Item(8, 15).Value2 = 0.4352965780925344
$ws.Cells.Item(8, 16).Value2 = 0.4352965780925344
$ws.Cells.Item(8, 17).Value2 = 91.02725908005468
$ws.Cells.Item(8, 18).Value2 = 819.245331720492
$ws.Cells.Item(8, 19).Value2 = 0.1152200080872524
$ws.Cells.Item(8, 20).Value2 = 0.1152200080872524
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 7).Value2 = 2.602283
$ws.Cells.Item(9, 8).Value2 = 7.806849
$ws.Cells.Item(9, 9).Value2 = 0.2646931170287289
$ws.Cells.Item(9, 10).Value2 = 0.2646931170287289
$ws.Cells.Item(9, 11).Value2 = 3
$ws.Cells.Item(9, 12).Value2 = 1
$ws.Cells.Item(9, 13).Value2 = 20.343383
$ws.Cells.Item(9, 14).Value2 = 61.03014900000001
$ws.Cells.Item(9, 15).Value2 = 0.2531579017099818
$ws.Cells.Item(9, 16).Value2 = 0.2531579017099818
$ws.Cells.Item(9, 17).Value2 = 52.939239743389
$ws.Cells.Item(9, 18).Value2 = 476.453157690501
$ws.Cells.Item(9, 19).Value2 = 0.06700915410406766
$ws.Cells.Item(9, 20).Value2 = 0.06700915410406766
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 7).Value2 = 2.602283
$ws.Cells.Item(10, 8).Value2 = 7.806849
$ws.Cells.Item(10, 9).Value2 = 0.2646931170287289
$ws.Cells.Item(10, 10).Value2 = 0.2646931170287289
$ws.Cells.Item(10, 11).Value2 = 3
$ws.Cells.Item(10, 12).Value2 = 1
$ws.Cells.Item(10, 13).Value2 = 25.035323
$ws.Cells.Item(10, 14).Value2 = 75.105969
$ws.Cells.Item(10, 15).Value2 = 0.3115455201974837
$ws.Cells.Item(10, 16).Value2 = 0.3115455201974837
$ws.Cells.Item(10, 17).Value2 = 65.148995442409
$ws.Cells.Item(10, 18).Value2 = 586.340958981681
$ws.Cells.Item(10, 19).Value2 = 0.08246395483740877
$ws.Cells.Item(10, 20).Value2 = 0.08246395483740877
$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 7).Value2 = 0.9246566666666666
$ws.Cells.Item(11, 8).Value2 = 2.77397
$ws.Cells.Item(11, 9).Value2 = 0.09405212856610688
$ws.Cells.Item(11, 10).Value2 = 0.0940521285661069
$ws.Cells.Item(11, 13).Value2 = 34.97976933333334
$ws.Cells.Item(11, 14).Value2 = 104.939308
$ws.Cells.Item(11, 15).Value2 = 0.4352965780925344
$ws.Cells.Item(11, 16).Value2 = 0.4352965780925344
$ws.Cells.Item(11, 17).Value2 = 32.34427691252889
$ws.Cells.Item(11, 18).Value2 = 291.09849221276
$ws.Cells.Item(11, 19).Value2 = 0.04094056972714543
$ws.Cells.Item(11, 20).Value2 = 0.04094056972714543
$ws.Cells.Item(12, 1).Value = "Resolving-Mac"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 7).Value2 = 0.9246566666666666
$ws.Cells.Item(12, 8).Value2 = 2.77397
$ws.Cells.Item(12, 9).Value2 = 0.09405212856610688
$ws.Cells.Item(12, 10).Value2 = 0.0940521285661069
$ws.Cells.Item(12, 13).Value2 = 20.343383
$ws.Cells.Item(12, 14).Value2 = 61.03014900000001
$ws.Cells.Item(12, 15).Value2 = 0.2531579017099818
$ws.Cells.Item(12, 16).Value2 = 0.2531579017099818
$ws.Cells.Item(12, 17).Value2 = 18.81064471350333
$ws.Cells.Item(12, 18).Value2 = 169.29580242153
$ws.Cells.Item(12, 19).Value2 = 0.02381003951915306
$ws.Cells.Item(12, 20).Value2 = 0.02381003951915306
$ws.Cells.Item(13, 1).Value = "Resolving-Mac"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 7).Value2 = 0.9246566666666666
$ws.Cells.Item(13, 8).Value2 = 2.77397
$ws.Cells.Item(13, 9).Value2 = 0.09405212856610688
$ws.Cells.Item(13, 10).Value2 = 0.0940521285661069
$ws.Cells.Item(13, 11).Value2 = 3
$ws.Cells.Item(13, 12).Value2 = 1
$ws.Cells.Item(13, 13).Value2 = 25.035323
$ws.Cells.Item(13, 14).Value2 = 75.105969
$ws.Cells.Item(13, 15).Value2 = 0.3115455201974837
$ws.Cells.Item(13, 16).Value2 = 0.3115455201974837
$ws.Cells.Item(13, 17).Value2 = 23.14907831410333
$ws.Cells.Item(13, 18).Value2 = 208.34170482693
$ws.Cells.Item(13, 19).Value2 = 0.02930151931980839
$ws.Cells.Item(13, 20).Value2 = 0.02930151931980839

# Remove now-obsolete rows (previously rows 14-17, target cluster "Resolving-Mac" no longer present)
$ws.Range("A14:T17").Delete()

